$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.005583666666666666
$ws.Range("H2").Value = 0.016751
$ws.Range("I2").Value = 0.0005349789730684028
$ws.Range("J2").Value = 0.0005349789730684027
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.502639
$ws.Range("N2").Value = 31.507917
$ws.Range("O2").Value = 0.9701256668284471
$ws.Range("P2").Value = 0.970125666828447
$ws.Range("Q2").Value = 0.05864323529633333
$ws.Range("R2").Value = 0.527789117667
$ws.Range("S2").Value = 0.0005189968329871821
$ws.Range("T2").Value = 0.000518996832987182

$ws.Range("G3").Value = 0.005583666666666666
$ws.Range("H3").Value = 0.016751
$ws.Range("I3").Value = 0.0005349789730684028
$ws.Range("J3").Value = 0.0005349789730684027
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2927206666666667
$ws.Range("N3").Value = 0.878162
$ws.Range("O3").Value = 0.02703852164627077
$ws.Range("P3").Value = 0.02703852164627077
$ws.Range("Q3").Value = 0.001634454629111111
$ws.Range("R3").Value = 0.014710091662
$ws.Range("S3").Value = 0.00001446504054360972
$ws.Range("T3").Value = 0.00001446504054360971

$ws.Range("G4").Value = 0.005583666666666666
$ws.Range("H4").Value = 0.016751
$ws.Range("I4").Value = 0.0005349789730684028
$ws.Range("J4").Value = 0.0005349789730684027
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03070066666666667
$ws.Range("N4").Value = 0.092102
$ws.Range("O4").Value = 0.002835811525282158
$ws.Range("P4").Value = 0.002835811525282158
$ws.Range("Q4").Value = 0.0001714222891111111
$ws.Range("R4").Value = 0.001542800602
$ws.Range("S4").Value = 0.00000151709953761099
$ws.Range("T4").Value = 0.00000151709953761099

$ws.Range("G5").Value = 7.96874
$ws.Range("H5").Value = 23.90622
$ws.Range("I5").Value = 0.7634962107066631
$ws.Range("J5").Value = 0.7634962107066631
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.502639
$ws.Range("N5").Value = 31.507917
$ws.Range("O5").Value = 0.9701256668284471
$ws.Range("P5").Value = 0.970125666828447
$ws.Range("Q5").Value = 83.69279950486001
$ws.Range("R5").Value = 753.23519554374
$ws.Range("S5").Value = 0.740687270532794
$ws.Range("T5").Value = 0.740687270532794

$ws.Range("G6").Value = 7.96874
$ws.Range("H6").Value = 23.90622
$ws.Range("I6").Value = 0.7634962107066631
$ws.Range("J6").Value = 0.7634962107066631
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2927206666666667
$ws.Range("N6").Value = 0.878162
$ws.Range("O6").Value = 0.02703852164627077
$ws.Range("P6").Value = 0.02703852164627077
$ws.Range("Q6").Value = 2.332614885293334
$ws.Range("R6").Value = 20.99353396764
$ws.Range("S6").Value = 0.02064380882003782
$ws.Range("T6").Value = 0.02064380882003782

$ws.Range("G7").Value = 7.96874
$ws.Range("H7").Value = 23.90622
$ws.Range("I7").Value = 0.7634962107066631
$ws.Range("J7").Value = 0.7634962107066631
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03070066666666667
$ws.Range("N7").Value = 0.092102
$ws.Range("O7").Value = 0.002835811525282158
$ws.Range("P7").Value = 0.002835811525282158
$ws.Range("Q7").Value = 0.2446456304933333
$ws.Range("R7").Value = 2.20181067444
$ws.Range("S7").Value = 0.00216513135383121
$ws.Range("T7").Value = 0.00216513135383121

$ws.Range("G8").Value = 2.462846666666667
$ws.Range("H8").Value = 7.388540000000001
$ws.Range("I8").Value = 0.2359688103202685
$ws.Range("J8").Value = 0.2359688103202685
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.502639
$ws.Range("N8").Value = 31.507917
$ws.Range("O8").Value = 0.9701256668284471
$ws.Range("P8").Value = 0.970125666828447
$ws.Range("Q8").Value = 25.86638945235334
$ws.Range("R8").Value = 232.79750507118
$ws.Range("S8").Value = 0.2289193994626658
$ws.Range("T8").Value = 0.2289193994626658

$ws.Range("G9").Value = 2.462846666666667
$ws.Range("H9").Value = 7.388540000000001
$ws.Range("I9").Value = 0.2359688103202685
$ws.Range("J9").Value = 0.2359688103202685
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2927206666666667
$ws.Range("N9").Value = 0.878162
$ws.Range("O9").Value = 0.02703852164627077
$ws.Range("P9").Value = 0.02703852164627077
$ws.Range("Q9").Value = 0.7209261181644446
$ws.Range("R9").Value = 6.488335063480001
$ws.Range("S9").Value = 0.006380247785689342
$ws.Range("T9").Value = 0.00638024778568934

$ws.Range("G10").Value = 2.462846666666667
$ws.Range("H10").Value = 7.388540000000001
$ws.Range("I10").Value = 0.2359688103202685
$ws.Range("J10").Value = 0.2359688103202685
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.03070066666666667
$ws.Range("N10").Value = 0.092102
$ws.Range("O10").Value = 0.002835811525282158
$ws.Range("P10").Value = 0.002835811525282158
$ws.Range("Q10").Value = 0.07561103456444446
$ws.Range("R10").Value = 0.6804993110800001
$ws.Range("S10").Value = 0.0006691630719133369
$ws.Range("T10").Value = 0.0006691630719133368
